$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.079.93'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.647.32'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("E4").Value = '  +0.57%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.81'
$ws.Range("E5").Value = '  +0.96%  '

$ws.Range("E6").Value = '  +0.89%  '

$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0640'
$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.67'
$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.876.31'
$ws.Range("E12").Value = '  +0.84%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.30'
$ws.Range("E13").Value = '  +1.64%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.673.55'
$ws.Range("E14").Value = '  +2.76%  '

$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  +1.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.49'
$ws.Range("E17").Value = '  +0.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.182.19'
$ws.Range("E18").Value = '  +1.01%  '

$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.51'
$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.82'
$ws.Range("E24").Value = '  +0.89%  '

$ws.Range("E25").Value = '  +3.84%  '

$ws.Range("E26").Value = '  +0.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '144.12'
$ws.Range("E27").Value = '  +0.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.91'
$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.58'
$ws.Range("E29").Value = '  +0.61%  '

$ws.Range("E30").Value = '  +1.35%  '

$ws.Range("E32").Value = '  +1.64%  '

$ws.Range("E33").Value = '  -0.29%  '

$ws.Range("E34").Value = '  -2.47%  '

$ws.Range("E35").Value = '  +1.16%  '

$ws.Range("E36").Value = '  +0.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.134.91'
$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("E38").Value = '  -1.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.47'

$ws.Range("E40").Value = '  +0.39%  '

$ws.Range("E41").Value = '  +1.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.48'
$ws.Range("E42").Value = '  +0.26%  '

$ws.Range("E43").Value = '  -0.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.785.77'
$ws.Range("E44").Value = '  +0.89%  '

$ws.Range("E45").Value = '  +4.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.75'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("E47").Value = '  +0.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.47'
$ws.Range("E48").Value = '  +0.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.78'
$ws.Range("E49").Value = '  +2.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.416'
$ws.Range("E50").Value = '  +0.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0962'
$ws.Range("E51").Value = '  +0.11%  '
